# feat: add S3 schedule
# Adds two new worksheets (S3_Sem_1, S3_Sem_2) at the end of the workbook,
# replicating the row/column layout used by every other sheet in the file
# (Hari, Jam, Ruangan, Prodi, Mata Kuliah, Semester, Kode Dosen, SKS).

$wb = $excel.ActiveWorkbook

$headerRow = @("Hari", "Jam", "Ruangan", "Prodi", "Mata Kuliah", "Semester", "Kode Dosen", "SKS")

$sheet17Data = @(
    @("SENIN", "10.00 - 11.50", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Topik Dalam Pengaman Jaringan B", "1", "TA", "3"),
    @("SENIN", "10.00 - 11.50", "Lab Pasca Lantai 1 (IF-110) - A`n(Kapasitas 8)", "S3", "Topik Dalam Data Mining A", "1", "AM", "3"),
    @("SELASA", "10.00 - 11.50", "Lab Pasca Lantai 1 (IF-110) - A`n(Kapasitas 8)", "S3", "Topik Dalam Sistem Terdistribusi A", "1", "AM", "3"),
    @("SELASA", "13.30 - 15.20", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Topik Dalam Data Mining C", "1", "CF, HF", "3"),
    @("RABU", "07.00 - 08.50", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Topik Dalam Data Deret Waktu A", "1", "BA", "3"),
    @("RABU", "07.00 - 08.50", "Lab Pasca Lantai 1 (IF-110) - A`n(Kapasitas 8)", "S3", "Topik Dalam Pengaman Jaringan A", "1", "HS, BJ", "3"),
    @("RABU", "10.00 - 11.50", "Lab Pasca Lantai 1 (IF-110) - A`n(Kapasitas 8)", "S3", "Topik Dalam Visi Komputer A", "1", "CF, WN", "3"),
    @("KAMIS", "07.00 - 08.50", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Filsafat Ilmu A", "1", "UY", "3"),
    @("KAMIS", "07.00 - 08.50", "Lab Pasca Lantai 1 (IF-110) - A`n(Kapasitas 8)", "S3", "Topik Dalam Text Mining A", "1", "DP", "3"),
    @("KAMIS", "10.00 - 11.50", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Topik Dalam Tata Kelola Teknologi Informasi T", "1", "RS", "3"),
    @("KAMIS", "10.00 - 11.50", "Lab Pasca Lantai 1 (IF-110) - A`n(Kapasitas 8)", "S3", "Topik Dalam Pengaman Jaringan B", "1", "TA", "3"),
    @("KAMIS", "13.30 - 15.20", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Topik Dalam Rekayasa Sistem Berbasis Pengetahuan T", "1", "RS", "3"),
    @("KAMIS", "09.00 - 10.50", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Topik Dalam Data Multivariat", "1", "BA", "3"),
    @("KAMIS", "09.00 - 10.50", "Lab Pasca Lantai 1 (IF-110) - A`n(Kapasitas 8)", "S3", "Topik Dalam Forensik Digital A", "1", "HS, BJ", "3"),
    @("KAMIS", "09.00 - 10.50", "Lab Pasca Lantai (IF-110) - B`n(Kapasitas 8)", "S3", "Topik Dalam Data Mining B", "1", "DP", "3")
)

$sheet18Data = @(
    @("RABU", "10.00 - 11.50", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Penulisan Ilmiah A", "2", "HS", "2"),
    @("RABU", "13.30 - 15.20", "IF-106 (kapasitas 33) / KHUSUS S3", "S3", "Metode Penelitian A", "2", "ST", "3")
)

function Add-ScheduleSheet {
    param(
        [string]$SheetName,
        [object[]]$HeaderRow,
        [object[]]$DataRows
    )

    $wb = $excel.ActiveWorkbook
    $afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
    $ws.Name = $SheetName

    $allRows = New-Object System.Collections.ArrayList
    [void]$allRows.Add($HeaderRow)
    foreach ($r in $DataRows) {
        [void]$allRows.Add($r)
    }

    for ($r = 0; $r -lt $allRows.Count; $r++) {
        $rowValues = $allRows[$r]
        for ($c = 0; $c -lt $rowValues.Count; $c++) {
            $val = $rowValues[$c]
            $cell = $ws.Cells.Item($r + 1, $c + 1)
            # The source data stores every value (including the purely numeric
            # looking ones, e.g. "1", "2", "3") as text/shared-strings rather
            # than numbers -- match that by forcing a text number format
            # before assigning anything that would otherwise be auto-coerced
            # into a number.
            if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $val
        }
    }

    # Drop the (now unused) text number format again so the sheet doesn't
    # carry any cell styling that wasn't in the original data.
    $ws.UsedRange.ClearFormats()
}

Add-ScheduleSheet "S3_Sem_1" $headerRow $sheet17Data
Add-ScheduleSheet "S3_Sem_2" $headerRow $sheet18Data
